# Update 'F' column (想去人数 / interest counts) values across all sheets
# per the refreshed data snapshot (commit: output generated at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 241
$ws.Range("F3").Value = 576
$ws.Range("F6").Value = 3170
$ws.Range("F7").Value = 2746
$ws.Range("F8").Value = 531
$ws.Range("F9").Value = 45
$ws.Range("F11").Value = 347
$ws.Range("F12").Value = 285
$ws.Range("F13").Value = 30
$ws.Range("F14").Value = 5679
$ws.Range("F15").Value = 617
$ws.Range("F16").Value = 1020
$ws.Range("F17").Value = 57
$ws.Range("F18").Value = 166
$ws.Range("F20").Value = 458
$ws.Range("F21").Value = 1234
$ws.Range("F22").Value = 75
$ws.Range("F23").Value = 6
$ws.Range("F24").Value = 121
$ws.Range("F25").Value = 333
$ws.Range("F26").Value = 44

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1169
$ws.Range("F4").Value = 14
$ws.Range("F5").Value = 23
$ws.Range("F7").Value = 23
$ws.Range("F9").Value = 57
$ws.Range("F21").Value = 52
$ws.Range("F23").Value = 339
$ws.Range("F29").Value = 208
$ws.Range("F30").Value = 60
$ws.Range("F34").Value = 16

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 82
$ws.Range("F4").Value = 43
$ws.Range("F5").Value = 2558
$ws.Range("F6").Value = 1113
$ws.Range("F9").Value = 1437
$ws.Range("F10").Value = 401
$ws.Range("F12").Value = 4

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 43
$ws.Range("F4").Value = 2558
$ws.Range("F6").Value = 1113
$ws.Range("F7").Value = 1437
$ws.Range("F8").Value = 401
$ws.Range("F10").Value = 14
$ws.Range("F11").Value = 241
$ws.Range("F12").Value = 576
$ws.Range("F13").Value = 23
$ws.Range("F14").Value = 3170
$ws.Range("F15").Value = 2746
$ws.Range("F16").Value = 531
$ws.Range("F17").Value = 45
$ws.Range("F20").Value = 23
$ws.Range("F21").Value = 347
$ws.Range("F23").Value = 57
$ws.Range("F24").Value = 30
$ws.Range("F25").Value = 5679
$ws.Range("F28").Value = 617
$ws.Range("F29").Value = 1020
$ws.Range("F31").Value = 57
$ws.Range("F32").Value = 166
$ws.Range("F38").Value = 52
$ws.Range("F39").Value = 339
$ws.Range("F40").Value = 1234
$ws.Range("F42").Value = 75
$ws.Range("F44").Value = 208
$ws.Range("F47").Value = 333
$ws.Range("F48").Value = 44
